$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column for UnitPrice between Item_Code (D) and HSN_Code (old E)
$ws.Columns("E:E").Insert()

# Header for new column E
$ws.Range("E1").Value = "UnitPrice"

# Apply number format (maps to built-in numFmtId 40) to the new UnitPrice column data cells
$ws.Range("E2:E3").NumberFormat = "#,##0.00_);[Red](#,##0.00)"

# Match the column width used for the rest of the sheet (closest attainable to 17.28515625)
$ws.Columns("E:E").ColumnWidth = 16.5

# Update row 2 data
$ws.Range("A2").Value = "Tomato454"
$ws.Range("B2").Value = "Buy"
$ws.Range("C2").Value = "Raw Material"
$ws.Range("D2").Value = "TM_0001434"
$ws.Range("E2").Value = 60.05
$ws.Range("F2").Value = 4512412
$ws.Range("G2").Value = 10000

# Update row 3 data
$ws.Range("A3").Value = "Ginger5786"
$ws.Range("B3").Value = "Sell"
$ws.Range("C3").Value = "Raw Material"
$ws.Range("D3").Value = "GIN_hjf"
$ws.Range("E3").Value = 70.78
$ws.Range("F3").Value = 451222
$ws.Range("G3").Value = 20000

# Fix up conditional formatting ranges (insert does not auto-shift these)
$fcs = $ws.Cells.FormatConditions
$fcs.Item(2).ModifyAppliesToRange($ws.Range("D2:E1048576"))
$fcs.Item(3).ModifyAppliesToRange($ws.Range("F2:F1048576"))

# Update selection to match target
[void]$ws.Range("A3").Select()
